$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header H1: "fess" -> "0010"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "0010"

# Columns H2:H31: change formula from SUBSTITUTE(...) to IF(YEAR(...)<2030, ...)
for ($row = 2; $row -le 31; $row++) {
    $ws.Range("H$row").Formula = "=IF(YEAR(G${row}:G${row})<2030, A${row}:A${row}, `"`")"
}
